# feat(CWL): add line number to build list
#
# The "General" sheet in this EN language-mod workbook lists CWL's
# localization keys (col A = id, col C = text_JP/source default, col D =
# text/English override). This change gives the BGM playlist control
# buttons and the converter "reload" button more descriptive English
# labels, and tweaks two warning/error message strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$nl = [char]10

# BGM playlist controls: generic one-word labels -> more descriptive ones
$ws.Range("D97").Value  = "View Playlist"   # was "View"
$ws.Range("D98").Value  = "Hide Playlist"   # was "Hide"
$ws.Range("D99").Value  = "Next Song"       # was "Next"
$ws.Range("D100").Value = "Last Song"       # was "Last"
$ws.Range("D102").Value = "Reload BGM"      # was "Reload"

# Converter panel reload button
$ws.Range("D104").Value = "Reload Data"     # was "Reload"

# Drama-play error message: drop "the" before "Player.log" (text_JP/C column)
$ws.Range("C112").Value = "Error occurred during drama play!" + $nl + "Please check Player.log and mods." + $nl + "{0}"

# Pop-empty-text warning: "prevented" -> "stopped"
$ws.Range("C114").Value = "{0} attempts to pop empty text, CWL stopped it"
$ws.Range("D114").Value = "{0} attempts to pop empty text, CWL stopped it"
